$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "Player Info" worksheet as the first sheet in the workbook
#    NOTE: worksheet references must be (re)fetched AFTER Worksheets.Add()
#    runs, otherwise stale references can end up aliasing the new sheet's
#    cells in this runtime.
# ---------------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"
$playerInfo.Move($wb.Worksheets.Item(1))

# Headers
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Copy the header style (bold, centered, bordered) from an existing sheet
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

$battingSheet.Range("A1").Copy()
$playerInfo.Range("A1:D1").PasteSpecial(-4122)

# Data row
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "4315"
$playerInfo.Range("A2").Style = "Normal"
$playerInfo.Range("B2").Value = "Nasir Jamal Ahmadzai"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Leg Break"

# ---------------------------------------------------------------------------
# 2. "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE, URLs -> bare codes
# ---------------------------------------------------------------------------
$battingSheet.Range("D1").Value = "MATCH_CODE"

$lastRow = $battingSheet.Cells.Item($battingSheet.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $battingSheet.Cells.Item($r, 4)
    $url = $cell.Value2
    if ($url) {
        $code = $url -replace ".*MatchCode=", ""
        $cell.NumberFormat = "@"
        $cell.Value = $code
        $cell.Style = "Normal"
    }
}

# ---------------------------------------------------------------------------
# 3. "ODI Bowling" sheet: MATCH_CARD_LINK -> MATCH_CODE, URLs -> bare codes
# ---------------------------------------------------------------------------
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$lastRowBowl = $bowlingSheet.Cells.Item($bowlingSheet.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRowBowl; $r++) {
    $cell = $bowlingSheet.Cells.Item($r, 2)
    $url = $cell.Value2
    if ($url) {
        $code = $url -replace ".*MatchCode=", ""
        $cell.NumberFormat = "@"
        $cell.Value = $code
        $cell.Style = "Normal"
    }
}

Write-Host "done"
